$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "MainMemId"
$ws.Range("C1").Value = "SubMemId"
$ws.Range("D1").Value = "RelateCode"
$ws.Range("E1").Value = "DateStart"

# Relationship rows (MainMemId / SubMemId / RelateCode)
$ws.Range("B2").Value = "0054042f-4d37-4800-a8ed-7ad25e23b764"
$ws.Range("C2").Value = "13b2f73c-65f2-4c1a-baa8-00e9cd5b18ac"
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = "13b2f73c-65f2-4c1a-baa8-00e9cd5b18ac"
$ws.Range("C3").Value = "0054042f-4d37-4800-a8ed-7ad25e23b764"
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = "0054042f-4d37-4800-a8ed-7ad25e23b764"
$ws.Range("C4").Value = "1a719c1b-8b56-42eb-adf2-4e40681d8806"
$ws.Range("D4").Value = 1

$ws.Range("B5").Value = "13b2f73c-65f2-4c1a-baa8-00e9cd5b18ac"
$ws.Range("C5").Value = "1a719c1b-8b56-42eb-adf2-4e40681d8806"
$ws.Range("D5").Value = 1

$ws.Range("B6").Value = "1a719c1b-8b56-42eb-adf2-4e40681d8806"
$ws.Range("C6").Value = "2d4bf196-a96a-47ed-a64a-f5f9c95be937"
$ws.Range("D6").Value = 0

$ws.Range("B7").Value = "2d4bf196-a96a-47ed-a64a-f5f9c95be937"
$ws.Range("C7").Value = "1a719c1b-8b56-42eb-adf2-4e40681d8806"
$ws.Range("D7").Value = 0

$ws.Range("B8").Value = "1a719c1b-8b56-42eb-adf2-4e40681d8806"
$ws.Range("C8").Value = "58d017d5-ef0c-4ea0-ae95-b8e50fa37b6d"
$ws.Range("D8").Value = 1

$ws.Range("B9").Value = "1a719c1b-8b56-42eb-adf2-4e40681d8806"
$ws.Range("C9").Value = "812527b5-2a2e-42bc-b5dd-9211bb750987"
$ws.Range("D9").Value = 1

$ws.Range("B10").Value = "1a719c1b-8b56-42eb-adf2-4e40681d8806"
$ws.Range("C10").Value = "c51a1ee3-aedf-4215-bb24-e23f332e4749"
$ws.Range("D10").Value = 1

$ws.Range("B11").Value = "1a719c1b-8b56-42eb-adf2-4e40681d8806"
$ws.Range("C11").Value = "d93088c1-9f47-4034-897f-10fdec8c5700"
$ws.Range("D11").Value = 1

$ws.Range("B12").Value = "2d4bf196-a96a-47ed-a64a-f5f9c95be937"
$ws.Range("C12").Value = "58d017d5-ef0c-4ea0-ae95-b8e50fa37b6d"
$ws.Range("D12").Value = 1

$ws.Range("B13").Value = "2d4bf196-a96a-47ed-a64a-f5f9c95be937"
$ws.Range("C13").Value = "812527b5-2a2e-42bc-b5dd-9211bb750987"
$ws.Range("D13").Value = 1

$ws.Range("B14").Value = "2d4bf196-a96a-47ed-a64a-f5f9c95be937"
$ws.Range("C14").Value = "c51a1ee3-aedf-4215-bb24-e23f332e4749"
$ws.Range("D14").Value = 1

$ws.Range("B15").Value = "2d4bf196-a96a-47ed-a64a-f5f9c95be937"
$ws.Range("C15").Value = "d93088c1-9f47-4034-897f-10fdec8c5700"
$ws.Range("D15").Value = 1

# Column widths (best-fit, matching widths Excel auto-sized for the UUID/text columns)
$ws.Columns("B:C").ColumnWidth = 36.451822916666664
$ws.Columns("D:D").ColumnWidth = 10.451822916666666
$ws.Columns("E:E").ColumnWidth = 8.451822916666666

# Restore the selection to match the saved view state
$ws.Range("G9").Select()
